$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.871.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.844.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "113.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  +1.81%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +4.16%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0851"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("E14").Value = "  +1.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.288.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.988"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.851.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.910.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  +9.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0975"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("E31").Value = "  +14.90%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0893"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.06%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.11%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.174.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +19.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.954"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.21%  "
